# Ambermoon-Advanced workbook update:
#   "Improved Tristan text, fixed chests in gatekeeper's house"
#
# 1) Todo sheet: the "Add gatekeeper chests (or fix them)" item is done -> remove it.
# 2) Maps sheet: add a new map entry for the Gatekeeper's House (3 houses in cavetown).

$wb = $excel.ActiveWorkbook

# --- Maps sheet: insert the new Gatekeeper's House map row ---
$maps = $wb.Worksheets.Item("Maps")
$maps.Rows.Item(13).Insert()
$maps.Cells.Item(13, 1).Value = 376
$maps.Cells.Item(13, 2).Value = "Pförtnerhaus / Gatekeeper's House"
$maps.Cells.Item(13, 3).Value = "2D"
$maps.Cells.Item(13, 4).Value = "3 houses in cavetown"
$maps.Range("A14").Select()

# --- Cosmetic selection updates matching the saved workbook view state ---
$chests = $wb.Worksheets.Item("Chests")
$chests.Range("E24").Select()

$places = $wb.Worksheets.Item("Places")
$places.Range("C8").Select()

# --- Todo sheet: remove the now-finished "Add gatekeeper chests" entry ---
# (selected/activated last so it remains the workbook's active sheet/tab,
# matching the saved file's tabSelected state)
$todo = $wb.Worksheets.Item("Todo")
$todo.Cells.Item(3, 1).ClearContents()
$todo.Range("A3").Select()
